$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rweText = "`nThe objective of this WG is to develop recommendations for generating robust RWE/RWD to support regulatory approval for CGT products. Specifically, the WG will review methods for direct and indirect comparisons using RWD and reflects on the opportunities and challenges of these approaches in the setting of CGT drug development using case studies as examples."
$trialDesignText = "The objective of this subgroup is to research and evaluate clinical trial design options in CGT development focusing on registrational trials. Specifically, we aim to explore new ideas and innovative approaches to provide insights, best practices, and recommendations to CGT drug development community on effective designs to accelerate drug approvals in various therapeutic areas including but not limited to oncology, autoimmune, and rare genetic disorders."

$ws.Range("B2").Value = $rweText
$ws.Range("B8").Value = $trialDesignText

$ws.Rows.Item(2).RowHeight = 60
$ws.Rows.Item(8).RowHeight = 60

$ws.Range("B6").Select()
